# Update the public EPEX Spot prices workbook with the latest day of data.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": add a new date column CP (15-sep) with hourly prices.
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Duplicate the style of the previous header cell (CO1) onto the new one (CP1)
# so the new header keeps the same bold / centered / bordered look.
$wsPrix.Range("CO1").Copy()
$wsPrix.Range("CP1").PasteSpecial(-4122)
$wsPrix.Range("CP1").Value = "15-sep"

$cpValues = @{
    2  = 0
    3  = -0.01
    4  = -0.01
    5  = -0.01
    6  = -0.01
    7  = -0.01
    8  = 0
    9  = 2.18
    10 = 19.69
    11 = 0.9
    12 = 0
    13 = -0.01
    14 = -0.04
    15 = -0.02
    16 = -0.06
    17 = -0.01
    18 = 0
    19 = 0.43
    20 = 6.03
    21 = 28.19
    22 = 9.369999999999999
    23 = 5.16
    24 = 5.17
    25 = 4.29
}

foreach ($r in $cpValues.Keys) {
    $wsPrix.Range("CP" + $r).Value = $cpValues[$r]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append two more daily rows (2025-09-13 and 2025-09-14), using
# the last known price (same as row 90 / 2025-09-12) carried forward.
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date column to stay as plain text (the sheet stores dates as
# text, not real Excel dates) and restore the default "Normal" style
# afterwards, just like the pre-existing rows above.
$wsGaz.Range("A91").NumberFormat = "@"
$wsGaz.Range("A91").Value = "2025-09-13"
$wsGaz.Range("A91").Style = "Normal"
$wsGaz.Range("B91").Value = 32.2

$wsGaz.Range("A92").NumberFormat = "@"
$wsGaz.Range("A92").Value = "2025-09-14"
$wsGaz.Range("A92").Style = "Normal"
$wsGaz.Range("B92").Value = 32.2

# ---------------------------------------------------------------------------
# Sheet "CO2": append two more daily rows (2025-09-13 and 2025-09-14), using
# the last known price (same as row 90 / 2025-09-12) carried forward.
# ---------------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")

$wsCO2.Range("A91").NumberFormat = "@"
$wsCO2.Range("A91").Value = "2025-09-13"
$wsCO2.Range("A91").Style = "Normal"
$wsCO2.Range("B91").Value = 75.47

$wsCO2.Range("A92").NumberFormat = "@"
$wsCO2.Range("A92").Value = "2025-09-14"
$wsCO2.Range("A92").Style = "Normal"
$wsCO2.Range("B92").Value = 75.47

Write-Output "Workbook updated: added CP column to Prix Spot and 2 rows to Gaz/CO2"
